# Scheduled-runner market-data refresh.
#
# Universalis price pull refreshed the cached currentAveragePrice(NQ/HQ) /
# LevePrice(NQ/HQ) / LeveProfit(NQ/HQ) figures (columns H-N) on each job's
# Leve-profit sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). These are
# plain cached numbers (no formulas feed them), so the refresh is just an
# in-place overwrite of the affected cells -- including a few rows whose
# M/N (profit) cell was previously absent/blank and now gets a value.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 11333.417
$ws.Range("I40").Value = 13777.889
$ws.Range("K40").Value = 13777.889
$ws.Range("M40").Value = -13602.889
$ws.Range("H64").Value = 55410.26
$ws.Range("J64").Value = 2873.1333
$ws.Range("L64").Value = 2873.1333
$ws.Range("N64").Value = -3369.1333
$ws.Range("H67").Value = 55410.26
$ws.Range("J67").Value = 2873.1333
$ws.Range("L67").Value = 2873.1333
$ws.Range("N67").Value = -4589.1333
$ws.Range("H132").Value = 18997.527
$ws.Range("I132").Value = 3314.4
$ws.Range("J132").Value = 89571.60000000001
$ws.Range("K132").Value = 9943.200000000001
$ws.Range("L132").Value = 268714.8
$ws.Range("M132").Value = -7413.200000000001
$ws.Range("N132").Value = -273774.8
$ws.Range("H136").Value = 53131.668
$ws.Range("I136").Value = 40000
$ws.Range("J136").Value = 59697.5
$ws.Range("K136").Value = 40000
$ws.Range("L136").Value = 59697.5
$ws.Range("M136").Value = -34900
$ws.Range("N136").Value = -69897.5
$ws.Range("H137").Value = 1322.9286
$ws.Range("I137").Value = 821.62067
$ws.Range("J137").Value = 1533.6232
$ws.Range("K137").Value = 2464.86201
$ws.Range("L137").Value = 4600.8696
$ws.Range("M137").Value = 85.13799000000017
$ws.Range("N137").Value = -9700.8696
$ws.Range("H138").Value = 2529.7397
$ws.Range("I138").Value = 2261.423
$ws.Range("J138").Value = 2678.1702
$ws.Range("K138").Value = 6784.268999999999
$ws.Range("L138").Value = 8034.5106
$ws.Range("M138").Value = -1644.268999999999
$ws.Range("N138").Value = -18314.5106

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 11281.667
$ws.Range("I28").Value = 2563.3333
$ws.Range("K28").Value = 2563.3333
$ws.Range("M28").Value = -2371.3333
$ws.Range("H99").Value = 11281.667
$ws.Range("I99").Value = 2563.3333
$ws.Range("K99").Value = 2563.3333
$ws.Range("M99").Value = 431.6667000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H98").Value = 20000
$ws.Range("J98").Value = 20000
$ws.Range("L98").Value = 20000
$ws.Range("N98").Value = -25990
$ws.Range("H107").Value = 2010.5385
$ws.Range("I107").Value = 1928.55
$ws.Range("J107").Value = 2283.8333
$ws.Range("K107").Value = 1928.55
$ws.Range("L107").Value = 2283.8333
$ws.Range("M107").Value = -8.549999999999955
$ws.Range("N107").Value = -6123.8333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 204329.25
$ws.Range("J4").Value = 226943.56
$ws.Range("L4").Value = 226943.56
$ws.Range("N4").Value = -227167.56
$ws.Range("H62").Value = 2667.5
$ws.Range("I62").Value = 2534.375
$ws.Range("J62").Value = 3200
$ws.Range("K62").Value = 2534.375
$ws.Range("L62").Value = 3200
$ws.Range("M62").Value = -1910.375
$ws.Range("N62").Value = -4448
$ws.Range("H65").Value = 2667.5
$ws.Range("I65").Value = 2534.375
$ws.Range("J65").Value = 3200
$ws.Range("K65").Value = 12671.875
$ws.Range("L65").Value = 16000
$ws.Range("M65").Value = -9551.875
$ws.Range("N65").Value = -22240
$ws.Range("H97").Value = 20000
$ws.Range("J97").Value = 20000
$ws.Range("L97").Value = 20000
$ws.Range("N97").Value = -21982
$ws.Range("H132").Value = 33616.387
$ws.Range("I132").Value = 1260.3143
$ws.Range("J132").Value = 159445.56
$ws.Range("K132").Value = 3780.9429
$ws.Range("L132").Value = 478336.68
$ws.Range("M132").Value = -1250.9429
$ws.Range("N132").Value = -483396.68

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1423.0526
$ws.Range("I5").Value = 766.0833
$ws.Range("K5").Value = 2298.2499
$ws.Range("M5").Value = -2186.2499
$ws.Range("H122").Value = 3124.7317
$ws.Range("I122").Value = 707.9231
$ws.Range("K122").Value = 6371.3079
$ws.Range("M122").Value = -3921.3079
$ws.Range("H131").Value = 908.4299999999999
$ws.Range("I131").Value = 711.8
$ws.Range("J131").Value = 918.7789299999999
$ws.Range("K131").Value = 2135.4
$ws.Range("L131").Value = 2756.33679
$ws.Range("M131").Value = 2904.6
$ws.Range("N131").Value = -12836.33679
$ws.Range("H135").Value = 1423.0526
$ws.Range("I135").Value = 766.0833
$ws.Range("K135").Value = 6894.7497
$ws.Range("M135").Value = -4359.7497

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 14109.167
$ws.Range("J5").Value = 14109.167
$ws.Range("L5").Value = 14109.167
$ws.Range("N5").Value = -14333.167
$ws.Range("H102").Value = 1475.8334
$ws.Range("I102").Value = 1130.375
$ws.Range("J102").Value = 1752.2
$ws.Range("K102").Value = 1130.375
$ws.Range("L102").Value = 1752.2
$ws.Range("M102").Value = 491.625
$ws.Range("N102").Value = -4996.2
$ws.Range("H113").Value = 2800
$ws.Range("J113").Value = 2800
$ws.Range("L113").Value = 2800
$ws.Range("N113").Value = -7140
$ws.Range("H126").Value = 6245.6924
$ws.Range("I126").Value = 10927
$ws.Range("J126").Value = 2233.1428
$ws.Range("K126").Value = 32781
$ws.Range("L126").Value = 6699.428400000001
$ws.Range("M126").Value = -30311
$ws.Range("N126").Value = -11639.4284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 6595.591
$ws.Range("J2").Value = 17857.285
$ws.Range("L2").Value = 17857.285
$ws.Range("N2").Value = -18081.285
$ws.Range("H7").Value = 3626.875
$ws.Range("I7").Value = 3602
$ws.Range("J7").Value = 3668.3333
$ws.Range("K7").Value = 3602
$ws.Range("L7").Value = 3668.3333
$ws.Range("M7").Value = -3490
$ws.Range("N7").Value = -3892.3333
$ws.Range("H40").Value = 2774.5
$ws.Range("I40").Value = 2659.4
$ws.Range("J40").Value = 2966.3333
$ws.Range("K40").Value = 2659.4
$ws.Range("L40").Value = 2966.3333
$ws.Range("M40").Value = -2523.4
$ws.Range("N40").Value = -3238.3333
$ws.Range("H61").Value = 3740
$ws.Range("I61").Value = 3675
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 3675
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -3473
$ws.Range("N61").Value = -4404
$ws.Range("H68").Value = 2984.6155
$ws.Range("I68").Value = 2450
$ws.Range("J68").Value = 3222.2222
$ws.Range("K68").Value = 2450
$ws.Range("L68").Value = 3222.2222
$ws.Range("M68").Value = -1701
$ws.Range("N68").Value = -4720.2222
$ws.Range("H71").Value = 2984.6155
$ws.Range("I71").Value = 2450
$ws.Range("J71").Value = 3222.2222
$ws.Range("K71").Value = 12250
$ws.Range("L71").Value = 16111.111
$ws.Range("M71").Value = -8506
$ws.Range("N71").Value = -23599.111
$ws.Range("H102").Value = 48557
$ws.Range("J102").Value = 48557
$ws.Range("L102").Value = 48557
$ws.Range("N102").Value = -55047
$ws.Range("H113").Value = 3740
$ws.Range("I113").Value = 3675
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 3675
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -1505
$ws.Range("N113").Value = -8340
$ws.Range("H126").Value = 3626.875
$ws.Range("I126").Value = 3602
$ws.Range("J126").Value = 3668.3333
$ws.Range("K126").Value = 10806
$ws.Range("L126").Value = 11004.9999
$ws.Range("M126").Value = -8336
$ws.Range("N126").Value = -15944.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2000
$ws.Range("J2").Value = 2000
$ws.Range("L2").Value = 2000
$ws.Range("N2").Value = -2224

Write-Output "Updated $($wb.Worksheets.Count) sheets: cached Leve market/profit figures refreshed."
